# Apply weekly update to the "Apio" (Celery) price sheet.
# A new week of data (one "Primera" + one "Segunda" quality row) is
# inserted at row 338, pushing all the existing historical rows down by
# two rows. The two oldest rows (formerly 389-390) end up appended as
# new rows 391-392 at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 338 (shifts 338:390 -> 340:392).
$ws.Rows.Item(338).Insert()
$ws.Rows.Item(338).Insert()

# --- Row 338: new "Primera" quality entry -----------------------------
$ws.Range("A338").Value = 6
$ws.Range("B338").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C338").Value = "Metropolitana"
$ws.Range("D338").Value = 44474
$ws.Range("E338").Value = 13
$ws.Range("F338").Value = 100112017
$ws.Range("G338").Value = "Apio"
$ws.Range("H338").Value = "Americana (o)"
$ws.Range("I338").Value = "Primera"
$ws.Range("J338").Value = 1700
$ws.Range("K338").Value = 6000
$ws.Range("L338").Value = 7000
$ws.Range("M338").Value = 6559
$ws.Range("N338").Value = "`$/docena de matas"
$ws.Range("O338").Value = "Región de Coquimbo"
$ws.Range("P338").Value = 1093
$ws.Range("Q338").Value = 6
$ws.Range("R338").Value = "Hortaliza"

# --- Row 339: new "Segunda" quality entry ------------------------------
$ws.Range("A339").Value = 6
$ws.Range("B339").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C339").Value = "Metropolitana"
$ws.Range("D339").Value = 44474
$ws.Range("E339").Value = 13
$ws.Range("F339").Value = 100112017
$ws.Range("G339").Value = "Apio"
$ws.Range("H339").Value = "Americana (o)"
$ws.Range("I339").Value = "Segunda"
$ws.Range("J339").Value = 550
$ws.Range("K339").Value = 5000
$ws.Range("L339").Value = 5000
$ws.Range("M339").Value = 5000
$ws.Range("N339").Value = "`$/docena de matas"
$ws.Range("O339").Value = "Región de Coquimbo"
$ws.Range("P339").Value = 833
$ws.Range("Q339").Value = 6
$ws.Range("R339").Value = "Hortaliza"
